$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.670.72'
$ws.Range("E2").Value = '  +1.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.868.95'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.75'
$ws.Range("E5").Value = '  +2.52%  '

$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4694'
$ws.Range("E7").Value = '  +3.99%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3940'
$ws.Range("E8").Value = '  +1.91%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.73'
$ws.Range("E9").Value = '  -0.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08058'
$ws.Range("E10").Value = '  +1.97%  '

$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.85'
$ws.Range("E12").Value = '  +1.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.875.87'
$ws.Range("E13").Value = '  +0.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.952'
$ws.Range("E14").Value = '  +1.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.144'
$ws.Range("E15").Value = '  -0.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001046'
$ws.Range("E17").Value = '  +1.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '86.63'
$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06620'
$ws.Range("E19").Value = '  +1.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.11'
$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.679.64'
$ws.Range("E22").Value = '  +1.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.503'
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.96'
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("E25").Value = '  +1.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.102.23'
$ws.Range("E26").Value = '  +0.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.75'
$ws.Range("E27").Value = '  +4.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.20'
$ws.Range("E28").Value = '  +2.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.092'
$ws.Range("E29").Value = '  +1.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.548'
$ws.Range("E30").Value = '  +0.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.35'
$ws.Range("E31").Value = '  +1.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9698'
$ws.Range("E32").Value = '  +3.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09505'
$ws.Range("E33").Value = '  +2.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.447'
$ws.Range("E34").Value = '  -3.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.588'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.328'
$ws.Range("E36").Value = '  +0.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02261'
$ws.Range("E37").Value = '  +1.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06096'
$ws.Range("E38").Value = '  +1.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.227'
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.137'
$ws.Range("E40").Value = '  -1.37%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5998'
$ws.Range("E41").Value = '  +1.56%  '

$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("E42").Value = '  +0.22%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1894'
$ws.Range("E43").Value = '  +0.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.23'
$ws.Range("E44").Value = '  +0.69%  '

$ws.Range("E45").Value = '  -0.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5680'
$ws.Range("E46").Value = '  +0.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.20'
$ws.Range("E47").Value = '  +2.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.395'
$ws.Range("E48").Value = '  +0.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.933'
$ws.Range("E49").Value = '  +0.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06849'
$ws.Range("E50").Value = '  +0.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '114.41'
$ws.Range("E51").Value = '  +5.92%  '
